$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.824.07"
$ws.Range("E2").Value = "  +1.64%  "
$ws.Range("D3").Value = "1.881.94"
$ws.Range("E3").Value = "  +1.77%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "333.66"
$ws.Range("E5").Value = "  +3.82%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "0.4749"
$ws.Range("E7").Value = "  +6.66%  "
$ws.Range("D8").Value = "0.3986"
$ws.Range("E8").Value = "  +4.32%  "
$ws.Range("D9").Value = "48.12"
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("D10").Value = "0.08063"
$ws.Range("E10").Value = "  +3.00%  "
$ws.Range("D11").Value = "1.027"
$ws.Range("E11").Value = "  +1.45%  "
$ws.Range("D12").Value = "21.94"
$ws.Range("E12").Value = "  +2.77%  "
$ws.Range("D13").Value = "1.893.43"
$ws.Range("E13").Value = "  +2.45%  "
$ws.Range("D14").Value = "5.972"
$ws.Range("E14").Value = "  +2.34%  "
$ws.Range("D15").Value = "7.192"
$ws.Range("E15").Value = "  +1.62%  "
$ws.Range("D16").Value = "1.004"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").Value = "0.00001055"
$ws.Range("E17").Value = "  +2.91%  "
$ws.Range("D18").Value = "87.23"
$ws.Range("E18").Value = "  +2.22%  "
$ws.Range("E19").Value = "  +1.86%  "
$ws.Range("D20").Value = "17.39"
$ws.Range("E20").Value = "  +2.80%  "
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").Value = "27.822.57"
$ws.Range("D23").Value = "5.512"
$ws.Range("E23").Value = "  +1.20%  "
$ws.Range("E24").Value = "  +2.73%  "
$ws.Range("D25").Value = "2.298"
$ws.Range("E25").Value = "  +1.97%  "
$ws.Range("D26").Value = "2.106.69"
$ws.Range("E26").Value = "  +1.93%  "
$ws.Range("D27").Value = "157.00"
$ws.Range("E27").Value = "  +3.33%  "
$ws.Range("D28").Value = "20.25"
$ws.Range("E28").Value = "  +4.93%  "
$ws.Range("D29").Value = "2.110"
$ws.Range("E29").Value = "  +3.21%  "
$ws.Range("D30").Value = "5.624"
$ws.Range("E30").Value = "  +2.08%  "
$ws.Range("D31").Value = "122.81"
$ws.Range("E31").Value = "  +2.94%  "
$ws.Range("D32").Value = "0.9774"
$ws.Range("E32").Value = "  +5.60%  "
$ws.Range("D33").Value = "0.09577"
$ws.Range("E33").Value = "  +2.91%  "
$ws.Range("D34").Value = "1.466"
$ws.Range("E34").Value = "  +1.37%  "
$ws.Range("D35").Value = "3.631"
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("D36").Value = "5.320"
$ws.Range("E36").Value = "  +1.84%  "
$ws.Range("E37").Value = "  +2.50%  "
$ws.Range("D38").Value = "0.06116"
$ws.Range("E38").Value = "  +3.04%  "
$ws.Range("D39").Value = "1.235"
$ws.Range("E39").Value = "  +3.07%  "
$ws.Range("D40").Value = "8.218"
$ws.Range("E40").Value = "  -0.75%  "
$ws.Range("D41").Value = "0.6032"
$ws.Range("E41").Value = "  +2.46%  "
$ws.Range("D42").Value = "1.001"
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").Value = "0.1911"
$ws.Range("E43").Value = "  +3.55%  "
$ws.Range("D44").Value = "10.32"
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "0.5745"
$ws.Range("E45").Value = "  +1.97%  "
$ws.Range("B46").Value = "WEMIXTOKEN"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "1.247"
$ws.Range("E46").Value = "  -0.85%  "
$ws.Range("D47").Value = "12.37"
$ws.Range("E47").Value = "  +2.17%  "
$ws.Range("D48").Value = "3.416"
$ws.Range("E48").Value = "  +1.78%  "
$ws.Range("D49").Value = "1.946"
$ws.Range("E49").Value = "  +1.41%  "
$ws.Range("D50").Value = "0.06810"
$ws.Range("E50").Value = "  -0.69%  "
$ws.Range("D51").Value = "113.60"
